# Populate the new "RunMode" column (C) on the LoginTestdata sheet with
# y/N values for the existing data rows, then leave the selection on B2
# (matching the saved view state of the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestdata")

$ws.Range("C2").Value = "y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "y"
$ws.Range("C5").Value = "y"

[void]$ws.Range("B2").Select()
